# Auto-generated edit script: update market-price snapshot values
# across ALC/ARM/BSM/CRP/CUL/GSM/LTW sheets (scheduled runner refresh).
$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 447.25
$ws.Range("I33").Value = 213.85
$ws.Range("K33").Value = 213.85
$ws.Range("M33").Value = 15.15000000000001
$ws.Range("H121").Value = 3500
$ws.Range("J121").Value = 3500
$ws.Range("L121").Value = 10500
$ws.Range("N121").Value = -13994
$ws.Range("H135").Value = 2462.9
$ws.Range("I135").Value = 2317.2727
$ws.Range("J135").Value = 2640.889
$ws.Range("K135").Value = 20855.4543
$ws.Range("L135").Value = 23768.001
$ws.Range("M135").Value = -18320.4543
$ws.Range("N135").Value = -28838.001
$ws.Range("H137").Value = 2249945.8
$ws.Range("I137").Value = 4047423
$ws.Range("J137").Value = 3099.35
$ws.Range("K137").Value = 12142269
$ws.Range("L137").Value = 9298.049999999999
$ws.Range("M137").Value = -12139719
$ws.Range("N137").Value = -14398.05
$ws.Range("H138").Value = 1634.15
$ws.Range("I138").Value = 881.09375
$ws.Range("J138").Value = 1988.5294
$ws.Range("K138").Value = 2643.28125
$ws.Range("L138").Value = 5965.5882
$ws.Range("M138").Value = 2496.71875
$ws.Range("N138").Value = -16245.5882
$ws.Range("H141").Value = 2688.3928
$ws.Range("I141").Value = 1983.0952
$ws.Range("K141").Value = 5949.2856
$ws.Range("M141").Value = -769.2856000000002

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 1074.2565
$ws.Range("I45").Value = 1092.8889
$ws.Range("J45").Value = 1032.3334
$ws.Range("K45").Value = 1092.8889
$ws.Range("L45").Value = 1032.3334
$ws.Range("M45").Value = -715.8888999999999
$ws.Range("N45").Value = -1786.3334
$ws.Range("H74").Value = 41248.42
$ws.Range("I74").Value = 45890.023
$ws.Range("J74").Value = 5662.8335
$ws.Range("K74").Value = 45890.023
$ws.Range("L74").Value = 5662.8335
$ws.Range("M74").Value = -45016.023
$ws.Range("N74").Value = -7410.8335
$ws.Range("H77").Value = 41248.42
$ws.Range("I77").Value = 45890.023
$ws.Range("J77").Value = 5662.8335
$ws.Range("K77").Value = 229450.115
$ws.Range("L77").Value = 28314.1675
$ws.Range("M77").Value = -225082.115
$ws.Range("N77").Value = -37050.1675

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H107").Value = 596.92
$ws.Range("I107").Value = 548.41174
$ws.Range("J107").Value = 700
$ws.Range("K107").Value = 548.41174
$ws.Range("L107").Value = 700
$ws.Range("M107").Value = 1371.58826
$ws.Range("N107").Value = -4540

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H16").Value = 2289.8667
$ws.Range("I16").Value = 2062.476
$ws.Range("J16").Value = 2820.4443
$ws.Range("K16").Value = 2062.476
$ws.Range("L16").Value = 2820.4443
$ws.Range("M16").Value = -1775.476
$ws.Range("N16").Value = -3394.4443
$ws.Range("H31").Value = 168460340
$ws.Range("I31").Value = 5003006
$ws.Range("J31").Value = 250189000
$ws.Range("K31").Value = 5003006
$ws.Range("L31").Value = 250189000
$ws.Range("M31").Value = -5002711
$ws.Range("N31").Value = -250189590
$ws.Range("H34").Value = 168460340
$ws.Range("I34").Value = 5003006
$ws.Range("J34").Value = 250189000
$ws.Range("K34").Value = 5003006
$ws.Range("L34").Value = 250189000
$ws.Range("M34").Value = -5002804
$ws.Range("N34").Value = -250189404
$ws.Range("H58").Value = 2551.88
$ws.Range("I58").Value = 2600.7727
$ws.Range("J58").Value = 2193.3333
$ws.Range("K58").Value = 2600.7727
$ws.Range("L58").Value = 2193.3333
$ws.Range("M58").Value = -2397.7727
$ws.Range("N58").Value = -2599.3333
$ws.Range("H107").Value = 4041.5366
$ws.Range("I107").Value = 361.4
$ws.Range("J107").Value = 5228.6772
$ws.Range("K107").Value = 361.4
$ws.Range("L107").Value = 5228.6772
$ws.Range("M107").Value = 1558.6
$ws.Range("N107").Value = -9068.6772
$ws.Range("H113").Value = 2289.8667
$ws.Range("I113").Value = 2062.476
$ws.Range("J113").Value = 2820.4443
$ws.Range("K113").Value = 2062.476
$ws.Range("L113").Value = 2820.4443
$ws.Range("M113").Value = 107.5239999999999
$ws.Range("N113").Value = -7160.4443
$ws.Range("H136").Value = 2551.88
$ws.Range("I136").Value = 2600.7727
$ws.Range("J136").Value = 2193.3333
$ws.Range("K136").Value = 7802.3181
$ws.Range("L136").Value = 6579.999899999999
$ws.Range("M136").Value = -5252.3181
$ws.Range("N136").Value = -11679.9999

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 3663480.2
$ws.Range("I4").Value = 6224816.5
$ws.Range("J4").Value = 4428.5713
$ws.Range("K4").Value = 18674449.5
$ws.Range("L4").Value = 13285.7139
$ws.Range("M4").Value = -18674337.5
$ws.Range("N4").Value = -13509.7139
$ws.Range("H36").Value = 0
$ws.Range("I36").Value = 0
$ws.Range("K36").Value = 0
$ws.Range("M36").Value = $null
$ws.Range("H122").Value = 1024.875
$ws.Range("I122").Value = 420
$ws.Range("J122").Value = 1111.2858
$ws.Range("K122").Value = 3780
$ws.Range("L122").Value = 10001.5722
$ws.Range("M122").Value = -1330
$ws.Range("N122").Value = -14901.5722
$ws.Range("H136").Value = 2107.1794
$ws.Range("I136").Value = 1518.4324
$ws.Range("K136").Value = 4555.2972
$ws.Range("M136").Value = 544.7028
$ws.Range("H137").Value = 22224760
$ws.Range("I137").Value = 1577.5555
$ws.Range("J137").Value = 55559536
$ws.Range("K137").Value = 4732.666499999999
$ws.Range("L137").Value = 166678608
$ws.Range("M137").Value = 367.3335000000006
$ws.Range("N137").Value = -166688808
$ws.Range("H138").Value = 20000
$ws.Range("I138").Value = 0
$ws.Range("K138").Value = 0
$ws.Range("M138").Value = $null

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H23").Value = 0
$ws.Range("I23").Value = 0
$ws.Range("J23").Value = 0
$ws.Range("K23").Value = 0
$ws.Range("L23").Value = 0
$ws.Range("M23").Value = $null
$ws.Range("N23").Value = $null
$ws.Range("H80").Value = 4631.8184
$ws.Range("I80").Value = 4631.8184
$ws.Range("J80").Value = 0
$ws.Range("K80").Value = 4631.8184
$ws.Range("L80").Value = 0
$ws.Range("M80").Value = -3633.8184
$ws.Range("N80").Value = $null
$ws.Range("H83").Value = 4631.8184
$ws.Range("I83").Value = 4631.8184
$ws.Range("J83").Value = 0
$ws.Range("K83").Value = 23159.092
$ws.Range("L83").Value = 0
$ws.Range("M83").Value = -18167.092
$ws.Range("N83").Value = $null
$ws.Range("H102").Value = 9955.429
$ws.Range("I102").Value = 9955.429
$ws.Range("K102").Value = 9955.429
$ws.Range("M102").Value = -8333.429
$ws.Range("H122").Value = 1855.0834
$ws.Range("I122").Value = 1625.2759
$ws.Range("J122").Value = 2807.1428
$ws.Range("K122").Value = 4875.8277
$ws.Range("L122").Value = 8421.428400000001
$ws.Range("M122").Value = -2425.8277
$ws.Range("N122").Value = -13321.4284
$ws.Range("H132").Value = 24920.164
$ws.Range("I132").Value = 28001.076
$ws.Range("J132").Value = 4509.125
$ws.Range("K132").Value = 84003.228
$ws.Range("L132").Value = 13527.375
$ws.Range("M132").Value = -81473.228
$ws.Range("N132").Value = -18587.375

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H82").Value = 2204.7
$ws.Range("I82").Value = 674.6667
$ws.Range("J82").Value = 4499.75
$ws.Range("K82").Value = 674.6667
$ws.Range("L82").Value = 4499.75
$ws.Range("M82").Value = -313.6667
$ws.Range("N82").Value = -5221.75
$ws.Range("H85").Value = 2204.7
$ws.Range("I85").Value = 674.6667
$ws.Range("J85").Value = 4499.75
$ws.Range("K85").Value = 674.6667
$ws.Range("L85").Value = 4499.75
$ws.Range("M85").Value = 573.3333
$ws.Range("N85").Value = -6995.75
$ws.Range("H93").Value = 71432550
$ws.Range("I93").Value = 606.375
$ws.Range("K93").Value = 606.375
$ws.Range("M93").Value = 641.625
$ws.Range("H95").Value = 53848
$ws.Range("J95").Value = 53848
$ws.Range("L95").Value = 53848
$ws.Range("N95").Value = -59340
$ws.Range("H132").Value = 6668.778
$ws.Range("I132").Value = 6007.4165
$ws.Range("J132").Value = 7197.8667
$ws.Range("K132").Value = 18022.2495
$ws.Range("L132").Value = 21593.6001
$ws.Range("M132").Value = -15492.2495
$ws.Range("N132").Value = -26653.6001

Write-Output "Applied scheduled runner price updates."